# Reorder the comma-separated "Recorded By" names in column G so that any
# token equal to "System" or "system" is moved to the front of the list,
# preserving the relative order of the remaining tokens.
#
# Example: "backup@backdoor.com, System, system" -> "System, system, backup@backdoor.com"
#          "dnasr281@gmail.com, System"           -> "System, dnasr281@gmail.com"
# Values that already start with System/system (or contain no System token)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ", "
    if ($parts.Count -lt 2) {
        continue
    }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p -eq "System" -or $p -eq "system") {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Count -eq 0) {
        continue
    }

    $newParts = $systemParts + $otherParts
    $newText = $newParts -join ", "

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
